{"js": "// Update the two-digit multiplication problems to the new values.\n// Each old expression appears exactly once in the document, so a\n// direct search + replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"59\u00d753=\", \"76\u00d765=\"],\n  [\"80\u00d744=\", \"34\u00d717=\"],\n  [\"59\u00d712=\", \"41\u00d712=\"],\n  [\"67\u00d744=\", \"95\u00d796=\"],\n  [\"27\u00d724=\", \"48\u00d740=\"],\n  [\"29\u00d761=\", \"85\u00d780=\"],\n  [\"84\u00d771=\", \"60\u00d722=\"],\n  [\"17\u00d718=\", \"21\u00d740=\"],\n  [\"48\u00d780=\", \"44\u00d736=\"],\n  [\"99\u00d792=\", \"82\u00d780=\"],\n  [\"98\u00d731=\", \"74\u00d743=\"],\n  [\"51\u00d752=\", \"52\u00d714=\"],\n  [\"11\u00d782=\", \"20\u00d718=\"],\n  [\"50\u00d731=\", \"20\u00d715=\"],\n  [\"59\u00d756=\", \"81\u00d787=\"],\n  [\"87\u00d715=\", \"18\u00d774=\"],\n  [\"47\u00d772=\", \"59\u00d772=\"],\n  [\"84\u00d745=\", \"39\u00d735=\"],\n  [\"58\u00d718=\", \"70\u00d780=\"],\n  [\"92\u00d771=\", \"81\u00d728=\"],\n  [\"60\u00d767=\", \"84\u00d786=\"],\n  [\"54\u00d718=\", \"75\u00d757=\"],\n  [\"76\u00d719=\", \"60\u00d788=\"],\n  [\"18\u00d738=\", \"74\u00d750=\"],\n  [\"77\u00d788=\", \"53\u00d768=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication problems to the new values.\n# Each old expression appears exactly once in the document, so a\n# direct Find/Replace per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n    \"59\u00d753=\" = \"76\u00d765=\";\n    \"80\u00d744=\" = \"34\u00d717=\";\n    \"59\u00d712=\" = \"41\u00d712=\";\n    \"67\u00d744=\" = \"95\u00d796=\";\n    \"27\u00d724=\" = \"48\u00d740=\";\n    \"29\u00d761=\" = \"85\u00d780=\";\n    \"84\u00d771=\" = \"60\u00d722=\";\n    \"17\u00d718=\" = \"21\u00d740=\";\n    \"48\u00d780=\" = \"44\u00d736=\";\n    \"99\u00d792=\" = \"82\u00d780=\";\n    \"98\u00d731=\" = \"74\u00d743=\";\n    \"51\u00d752=\" = \"52\u00d714=\";\n    \"11\u00d782=\" = \"20\u00d718=\";\n    \"50\u00d731=\" = \"20\u00d715=\";\n    \"59\u00d756=\" = \"81\u00d787=\";\n    \"87\u00d715=\" = \"18\u00d774=\";\n    \"47\u00d772=\" = \"59\u00d772=\";\n    \"84\u00d745=\" = \"39\u00d735=\";\n    \"58\u00d718=\" = \"70\u00d780=\";\n    \"92\u00d771=\" = \"81\u00d728=\";\n    \"60\u00d767=\" = \"84\u00d786=\";\n    \"54\u00d718=\" = \"75\u00d757=\";\n    \"76\u00d719=\" = \"60\u00d788=\";\n    \"18\u00d738=\" = \"74\u00d750=\";\n    \"77\u00d788=\" = \"53\u00d768=\";\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
